$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert the three new price-report rows.
#   - a single new row is inserted at row 126 (pushing the former 126..170
#     block down by one, to 127..171)
#   - two more new rows are inserted at 134-135 (pushing the former
#     133..170-block-already-shifted-by-one further down by two, ending at
#     136..173)
# ---------------------------------------------------------------------------
$ws.Rows.Item(126).Insert()
$ws.Rows.Item(134).Resize(2, 1).EntireRow.Insert()

function Set-GrapeRow($Row, $Date, $Variedad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidad) {
    $ws.Cells.Item($Row, 1).Value = 2
    $ws.Cells.Item($Row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($Row, 3).Value = "Coquimbo"
    $ws.Cells.Item($Row, 4).Value = $Date
    $ws.Cells.Item($Row, 5).Value = 4
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100109
    $ws.Cells.Item($Row, 8).Value = "Uva"
    $ws.Cells.Item($Row, 9).Value = 100109001
    $ws.Cells.Item($Row, 10).Value = "Uva"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = "Primera"
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# New row 126: 2022-12-22, Superior Seedless
Set-GrapeRow 126 44917 "Superior Seedless" 340 13500 14000 13750 "$/bandeja 10 kilos" "Provincia de Limarí" 1375 10

# New row 134: 2022-12-29, Flame Seedless
Set-GrapeRow 134 44924 "Flame Seedless" 1100 7500 8000 7750 "$/bandeja 10 kilos" "Provincia de Limarí" 775 10

# New row 135: 2022-12-29, Superior Seedless (Provincia de Huasco)
Set-GrapeRow 135 44924 "Superior Seedless" 1400 10500 11000 10750 "$/bandeja 10 kilos" "Provincia de Huasco" 1075 10
